# Adds an "Evaluation Warning" sheet (the classic Aspose.Cells unlicensed-use
# watermark sheet) after the existing "Empty" sheet, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet right after the last existing sheet -----------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Evaluation Warning"

# --- 2. Write the evaluation-warning text into A5 ----------------------------
$cell = $ws.Range("A5")
$cell.Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2020 Aspose Pty Ltd."

# --- 3. Format the cell: bold, italic, 18pt, blue ----------------------------
$cell.Font.Bold = $true
$cell.Font.Italic = $true
$cell.Font.Size = 18
$cell.Font.Color = 16711680

# Touch alignment/protection at their (already default) values so the xf
# records the same "apply*" bookkeeping flags the source file carries.
$cell.HorizontalAlignment = 1
$cell.Locked = $true

# --- 4. Row height for the banner row ----------------------------------------
$ws.Rows.Item(5).RowHeight = 23.25

# --- 5. Page setup: margins + portrait orientation ---------------------------
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
$ws.PageSetup.Orientation = 1

Write-Host "Added 'Evaluation Warning' sheet"
